$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.079.97'
$ws.Range("E2").Value = '  -0.42%  '

# Row 3
$ws.Range("D3").Value = '3.551.93'
$ws.Range("E3").Value = '  -1.59%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = "'587.66"
$ws.Range("E5").Value = '  +1.23%  '

# Row 6
$ws.Range("D6").Value = "'185.97"
$ws.Range("E6").Value = '  -2.12%  '

# Row 7
$ws.Range("D7").Value = '3.541.26'
$ws.Range("E7").Value = '  -1.78%  '

# Row 8
$ws.Range("D8").Value = "'0.616"
$ws.Range("E8").Value = '  -2.18%  '

# Row 9
$ws.Range("E9").Value = '  -0.08%  '

# Row 10
$ws.Range("E10").Value = '  +5.58%  '

# Row 11
$ws.Range("D11").Value = "'0.646"
$ws.Range("E11").Value = '  -1.90%  '

# Row 12
$ws.Range("D12").Value = "'54.44"
$ws.Range("E12").Value = '  -2.93%  '

# Row 13
$ws.Range("E13").Value = '  -2.01%  '

# Row 14
$ws.Range("D14").Value = "'9.50"
$ws.Range("E14").Value = '  -2.24%  '

# Row 15
$ws.Range("D15").Value = '4.109.79'
$ws.Range("E15").Value = '  -1.93%  '

# Row 16
$ws.Range("D16").Value = "'19.42"
$ws.Range("E16").Value = '  -2.08%  '

# Row 17
$ws.Range("D17").Value = '70.026.33'
$ws.Range("E17").Value = '  -0.54%  '

# Row 18
$ws.Range("D18").Value = '3.542.11'
$ws.Range("E18").Value = '  -1.85%  '

# Row 19
$ws.Range("D19").Value = "'12.48"
$ws.Range("E19").Value = '  -1.79%  '

# Row 20
$ws.Range("E20").Value = '  -1.28%  '

# Row 21
$ws.Range("D21").Value = "'539.99"
$ws.Range("E21").Value = '  +9.22%  '

# Row 22
$ws.Range("D22").Value = "'1.02"
$ws.Range("E22").Value = '  -2.78%  '

# Row 23
$ws.Range("D23").Value = "'18.01"
$ws.Range("E23").Value = '  -6.41%  '

# Row 24
$ws.Range("D24").Value = "'4.62"
$ws.Range("E24").Value = '  +5.48%  '

# Row 25
$ws.Range("E25").Value = '  -1.41%  '

# Row 26
$ws.Range("D26").Value = "'95.51"
$ws.Range("E26").Value = '  -0.82%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = "'11.26"
$ws.Range("E27").Value = '  +1.35%  '

# Row 28
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = "'3.00"
$ws.Range("E28").Value = '  +0.08%  '

# Row 29
$ws.Range("E29").Value = '  -3.72%  '

# Row 30
$ws.Range("D30").Value = "'32.16"
$ws.Range("E30").Value = '  -0.22%  '

# Row 31
$ws.Range("D31").Value = "'7.31"
$ws.Range("E31").Value = '  -3.83%  '

# Row 32
$ws.Range("D32").Value = "'12.48"
$ws.Range("E32").Value = '  +2.17%  '

# Row 33
$ws.Range("D33").Value = "'64.53"
$ws.Range("E33").Value = '  -1.88%  '

# Row 34
$ws.Range("E34").Value = '  -3.80%  '

# Row 35
$ws.Range("E35").Value = '  +7.70%  '

# Row 36
$ws.Range("D36").Value = "'549.73"
$ws.Range("E36").Value = '  -4.30%  '

# Row 37
$ws.Range("D37").Value = "'0.415"
$ws.Range("E37").Value = '  +4.27%  '

# Row 38
$ws.Range("D38").Value = "'38.47"
$ws.Range("E38").Value = '  -0.29%  '

# Row 39
$ws.Range("E39").Value = '  -0.18%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0767'
$ws.Range("E40").Value = '  -5.85%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = "'3.41"
$ws.Range("E41").Value = '  -3.82%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = "'0.135"
$ws.Range("E42").Value = '  -1.99%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = "'3.10"
$ws.Range("E43").Value = '  -8.19%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.345.13'
$ws.Range("E44").Value = '  +3.36%  '

# Row 45
$ws.Range("D45").Value = "'3.61"
$ws.Range("E45").Value = '  +6.39%  '

# Row 46
$ws.Range("E46").Value = '  -2.40%  '

# Row 47
$ws.Range("D47").Value = "'0.0442"
$ws.Range("E47").Value = '  -0.51%  '

# Row 48
$ws.Range("D48").Value = "'9.19"
$ws.Range("E48").Value = '  -6.39%  '

# Row 49
$ws.Range("E49").Value = '  -1.53%  '

# Row 50
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = '  +0.00%  '

# Row 51
$ws.Range("D51").Value = "'137.28"
$ws.Range("E51").Value = '  +1.37%  '
